$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q1" positioned right before "总计"
#    (i.e. right after "2021-Q4", which is the 3rd sheet).
#    NOTE: worksheet variables captured via Worksheets.Item(N) track
#    position N, not object identity - so we must (re)fetch sheet handles
#    AFTER the sheet collection size has changed, and only cache them
#    once no further Add()/Delete() calls will happen.
# ---------------------------------------------------------------------------
$q4_2021 = $wb.Worksheets.Item(3)
$newSheet = $wb.Worksheets.Add($null, $q4_2021)
$newSheet.Name = "2022-Q1"

# From this point on the sheet collection size is stable, so positional
# handles are safe to keep around.
$q4_2021 = $wb.Worksheets.Item(3)      # "2021-Q4"  (style/template source)
$newSheet = $wb.Worksheets.Item(4)     # "2022-Q1"  (just created)
$totalSheet = $wb.Worksheets.Item(5)   # "总计"

# Copy the look & feel (fonts / borders / alignment) of the existing
# per-quarter sheets so the new sheet matches the established style.
$q4_2021.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

$q4_2021.Range("A2:A4").Copy()
$newSheet.Range("A2:A10").PasteSpecial(-4122)  # xlPasteFormats

# ---------------------------------------------------------------------------
# 2. Fill in the header row
# ---------------------------------------------------------------------------
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# ---------------------------------------------------------------------------
# 3. Fill in the fund holdings data (rows 2-10)
#    Columns B-G are stored as text (matching the other quarter sheets),
#    column A is the 0-based row index and column H is a plain number.
# ---------------------------------------------------------------------------
$rows = @(
    @("160106", "南方高增长混合(LOF)", "20.27", "87.01", "3.30", "0.6689", 10),
    @("004616", "中欧电子信息产业沪港深股票A", "14.54", "92.26", "3.14", "0.4566", 8),
    @("160642", "鹏华增瑞灵活配置混合(LOF)", "6.76", "91.34", "5.07", "0.3427", 6),
    @("160105", "南方积极配置混合(LOF)", "7.84", "88.09", "3.34", "0.2619", 10),
    @("005763", "中欧电子信息产业沪港深股票C", "7.73", "92.26", "3.14", "0.2427", 8),
    @("159851", "华宝中证金融科技主题ETF", "3.16", "98.58", "2.82", "0.0891", 10),
    @("000554", "南方中国梦灵活配置混合", "1.39", "88.18", "3.25", "0.0452", 10),
    @("930602", "国信价值智选混合型集合资产管理计划", "0.50", "67.38", "4.94", "0.0247", 5),
    @("516100", "华夏中证金融科技主题交易型开放式指数证券投资基金", "0.68", "96.91", "2.80", "0.0190", 10)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $newSheet.Range("A$r").Value = $i
    $newSheet.Range("B$r").Value = "'" + $data[0]
    $newSheet.Range("B$r").Style = "Normal"
    $newSheet.Range("C$r").Value = $data[1]
    $newSheet.Range("D$r").Value = "'" + $data[2]
    $newSheet.Range("D$r").Style = "Normal"
    $newSheet.Range("E$r").Value = "'" + $data[3]
    $newSheet.Range("E$r").Style = "Normal"
    $newSheet.Range("F$r").Value = "'" + $data[4]
    $newSheet.Range("F$r").Style = "Normal"
    $newSheet.Range("G$r").Value = "'" + $data[5]
    $newSheet.Range("G$r").Style = "Normal"
    $newSheet.Range("H$r").Value = $data[6]
}

# Re-apply the column-A index formatting (the "Style = Normal" calls above do
# not touch column A, but let's make sure it keeps the bold/bordered look
# after all the edits above).
$q4_2021.Range("A2:A4").Copy()
$newSheet.Range("A2:A10").PasteSpecial(-4122)  # xlPasteFormats
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $newSheet.Range("A$r").Value = $i
}

# ---------------------------------------------------------------------------
# 4. Update the "总计" (totals) sheet: insert a new row for 2022-Q1 at the
#    top of the data (row 2), pushing the existing rows down.
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert(-4121)  # xlShiftDown

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$totalSheet.Range("B3:D3").Copy()
$totalSheet.Range("B2:D2").PasteSpecial(-4122)  # xlPasteFormats

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 9
$totalSheet.Range("D2").Value = 2.15

# Renumber the 0-based index column for the rows that shifted down.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3

Write-Host "Edit complete"
